$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# New columns were inserted (asset_description, location, amc_start_date,
# amc_end_date, warranty_start_date, warranty_end_date) and the existing
# header labels shifted to their new positions.
$ws.Range("A1").Value = "serial_number"
$ws.Range("B1").Value = "asset_type"
$ws.Range("C1").Value = "asset_description"
$ws.Range("D1").Value = "po_number"
$ws.Range("E1").Value = "sap_asset_id"
$ws.Range("F1").Value = "installation_date"
$ws.Range("G1").Value = "location"
$ws.Range("H1").Value = "amc_start_date"
$ws.Range("I1").Value = "amc_end_date"
$ws.Range("J1").Value = "warranty_start_date"
$ws.Range("K1").Value = "warranty_end_date"
$ws.Range("L1").Value = "amc_contract"
$ws.Range("M1").Value = "end_user"

# The original header style (bold font + border) only covered A1:G1 - extend
# it across the newly inserted header cells (H1:M1) by copying the format.
$ws.Range("A1").Copy()
$ws.Range("H1:M1").PasteSpecial(-4122)

# ---- Clear the old data row (row 2 previously held a single sample row) ----
$ws.Range("A2:G2").ClearContents()

# ---- Data rows ----
# Row 2 - Laptop
$ws.Range("A2").Value = "SN12345"
$ws.Range("B2").Value = "Laptop"
$ws.Range("C2").Value = "Dell Latitude 5520"
$ws.Range("D2").Value = 7912775
$ws.Range("E2").Value = "SAP123"
$ws.Range("H2").Value = (Get-Date -Year 2023 -Month 1 -Day 3 -Hour 0 -Minute 0 -Second 0).Date

# Row 3 - Server
$ws.Range("A3").Value = "SN67890"
$ws.Range("B3").Value = "Server"
$ws.Range("C3").Value = "HPE ProLiant DL360"
$ws.Range("D3").Value = 7050090689
$ws.Range("E3").Value = "SAP456"
$ws.Range("F3").Value = (Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("G3").Value = "EnD"
$ws.Range("M3").Value = "Tituraj"

# Row 4 - Printer
$ws.Range("A4").Value = "SN11223"
$ws.Range("B4").Value = "Printer"
$ws.Range("C4").Value = "Cisco ISR 4331"
$ws.Range("E4").Value = "SAP789"
$ws.Range("F4").Value = (Get-Date -Year 2021 -Month 2 -Day 2 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("G4").Value = "Shastri Bhawan"

# ---- Column widths for the newly populated columns ----
$ws.Range("C1:H1").EntireColumn.AutoFit()

# ---- Restore the active selection ----
$ws.Range("F5").Select()
